# Refresh cached market-price / profit figures on the Leve profit-tracking
# sheets (one per crafting job) to match the latest Universalis price pull.
# Columns: H/I/J = currentAveragePrice(/NQ/HQ), K/L = LevePrice(NQ/HQ),
#          M/N = LeveProfit(NQ/HQ).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 261194.17
$ws.Range("J17").Value = 281170.66
$ws.Range("L17").Value = 843511.98
$ws.Range("N17").Value = -843847.98
$ws.Range("H19").Value = 6999.4287
$ws.Range("I19").Value = 9399.200000000001
$ws.Range("K19").Value = 9399.200000000001
$ws.Range("M19").Value = -9224.200000000001
$ws.Range("H33").Value = 667.0833
$ws.Range("I33").Value = 733.8889
$ws.Range("J33").Value = 466.66666
$ws.Range("K33").Value = 733.8889
$ws.Range("L33").Value = 466.66666
$ws.Range("M33").Value = -504.8889
$ws.Range("N33").Value = -924.66666
$ws.Range("H54").Value = 23397.8
$ws.Range("I54").Value = 23397.8
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 23397.8
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -22911.8
$ws.Range("N54").ClearContents()
$ws.Range("H74").Value = 4111.8335
$ws.Range("I74").Value = 3223.6667
$ws.Range("K74").Value = 3223.6667
$ws.Range("M74").Value = -2287.6667
$ws.Range("H77").Value = 4111.8335
$ws.Range("I77").Value = 3223.6667
$ws.Range("K77").Value = 16118.3335
$ws.Range("M77").Value = -11438.3335
$ws.Range("H137").Value = 2577
$ws.Range("I137").Value = 2568.5
$ws.Range("J137").Value = 2586.4443
$ws.Range("K137").Value = 7705.5
$ws.Range("L137").Value = 7759.3329
$ws.Range("M137").Value = -5155.5
$ws.Range("N137").Value = -12859.3329
$ws.Range("H138").Value = 1480.2738
$ws.Range("I138").Value = 1266.5
$ws.Range("J138").Value = 1827.6562
$ws.Range("K138").Value = 3799.5
$ws.Range("L138").Value = 5482.9686
$ws.Range("M138").Value = 1340.5
$ws.Range("N138").Value = -15762.9686
$ws.Range("H141").Value = 21534.166
$ws.Range("I141").Value = 21720.457
$ws.Range("K141").Value = 65161.371
$ws.Range("M141").Value = -59981.371

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1365.5714
$ws.Range("I61").Value = 1365.5714
$ws.Range("K61").Value = 1365.5714
$ws.Range("M61").Value = -1153.5714
$ws.Range("H74").Value = 2281.6428
$ws.Range("I74").Value = 2078.5833
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 2078.5833
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1204.5833
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 2281.6428
$ws.Range("I77").Value = 2078.5833
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 10392.9165
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -6024.916499999999
$ws.Range("N77").Value = -26236
$ws.Range("H88").Value = 885.95654
$ws.Range("I88").Value = 588.9
$ws.Range("J88").Value = 1114.4615
$ws.Range("K88").Value = 588.9
$ws.Range("L88").Value = 1114.4615
$ws.Range("M88").Value = -182.9
$ws.Range("N88").Value = -1926.4615
$ws.Range("H91").Value = 885.95654
$ws.Range("I91").Value = 588.9
$ws.Range("J91").Value = 1114.4615
$ws.Range("K91").Value = 588.9
$ws.Range("L91").Value = 1114.4615
$ws.Range("M91").Value = 815.1
$ws.Range("N91").Value = -3922.4615
$ws.Range("H122").Value = 2000.5
$ws.Range("I122").Value = 1931.75
$ws.Range("K122").Value = 5795.25
$ws.Range("M122").Value = -3345.25
$ws.Range("H136").Value = 1365.5714
$ws.Range("I136").Value = 1365.5714
$ws.Range("K136").Value = 4096.7142
$ws.Range("M136").Value = -1546.7142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2382.4167
$ws.Range("I86").Value = 2386.75
$ws.Range("K86").Value = 2386.75
$ws.Range("M86").Value = -1263.75
$ws.Range("H89").Value = 2382.4167
$ws.Range("I89").Value = 2386.75
$ws.Range("K89").Value = 11933.75
$ws.Range("M89").Value = -6317.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32828.8
$ws.Range("I31").Value = 36409.1
$ws.Range("K31").Value = 36409.1
$ws.Range("M31").Value = -36114.1
$ws.Range("H34").Value = 32828.8
$ws.Range("I34").Value = 36409.1
$ws.Range("K34").Value = 36409.1
$ws.Range("M34").Value = -36207.1
$ws.Range("H51").Value = 13899
$ws.Range("J51").Value = 13899
$ws.Range("L51").Value = 13899
$ws.Range("N51").Value = -15371
$ws.Range("H58").Value = 491.69565
$ws.Range("I58").Value = 497
$ws.Range("J58").Value = 375
$ws.Range("K58").Value = 497
$ws.Range("L58").Value = 375
$ws.Range("M58").Value = -294
$ws.Range("N58").Value = -781
$ws.Range("H61").Value = 13899
$ws.Range("J61").Value = 13899
$ws.Range("L61").Value = 13899
$ws.Range("N61").Value = -14595
$ws.Range("H99").Value = 26972.54
$ws.Range("I99").Value = 30285.727
$ws.Range("J99").Value = 8750
$ws.Range("K99").Value = 30285.727
$ws.Range("L99").Value = 8750
$ws.Range("M99").Value = -28787.727
$ws.Range("N99").Value = -11746
$ws.Range("H105").Value = 1706.125
$ws.Range("I105").Value = 1446
$ws.Range("K105").Value = 1446
$ws.Range("M105").Value = 301
$ws.Range("H122").Value = 88291.336
$ws.Range("I122").Value = 130543.25
$ws.Range("J122").Value = 3787.5
$ws.Range("K122").Value = 391629.75
$ws.Range("L122").Value = 11362.5
$ws.Range("M122").Value = -389179.75
$ws.Range("N122").Value = -16262.5
$ws.Range("H126").Value = 26972.54
$ws.Range("I126").Value = 30285.727
$ws.Range("J126").Value = 8750
$ws.Range("K126").Value = 90857.181
$ws.Range("L126").Value = 26250
$ws.Range("M126").Value = -88387.181
$ws.Range("N126").Value = -31190
$ws.Range("H132").Value = 1871.5555
$ws.Range("I132").Value = 1995.3429
$ws.Range("K132").Value = 5986.028700000001
$ws.Range("M132").Value = -3456.028700000001
$ws.Range("H134").Value = 2060.5671
$ws.Range("I134").Value = 1937.875
$ws.Range("J134").Value = 2685.182
$ws.Range("K134").Value = 5813.625
$ws.Range("L134").Value = 8055.545999999999
$ws.Range("M134").Value = -3278.625
$ws.Range("N134").Value = -13125.546
$ws.Range("H136").Value = 491.69565
$ws.Range("I136").Value = 497
$ws.Range("J136").Value = 375
$ws.Range("K136").Value = 1491
$ws.Range("L136").Value = 1125
$ws.Range("M136").Value = 1059
$ws.Range("N136").Value = -6225
$ws.Range("H141").Value = 139152.86
$ws.Range("J141").Value = 149845.5
$ws.Range("L141").Value = 149845.5
$ws.Range("N141").Value = -160205.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 70692.55499999999
$ws.Range("I129").Value = 118128.06
$ws.Range("J129").Value = 3492.25
$ws.Range("K129").Value = 354384.18
$ws.Range("L129").Value = 10476.75
$ws.Range("M129").Value = -349384.18
$ws.Range("N129").Value = -20476.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26640
$ws.Range("H70").Value = 4034.4
$ws.Range("I70").Value = 3940.4736
$ws.Range("J70").Value = 4331.8335
$ws.Range("K70").Value = 3940.4736
$ws.Range("L70").Value = 4331.8335
$ws.Range("M70").Value = -3670.4736
$ws.Range("N70").Value = -4871.8335
$ws.Range("H73").Value = 4034.4
$ws.Range("I73").Value = 3940.4736
$ws.Range("J73").Value = 4331.8335
$ws.Range("K73").Value = 3940.4736
$ws.Range("L73").Value = 4331.8335
$ws.Range("M73").Value = -3004.4736
$ws.Range("N73").Value = -6203.8335
$ws.Range("H120").Value = 90000
$ws.Range("I120").Value = 90000
$ws.Range("K120").Value = 90000
$ws.Range("M120").Value = -85162

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 830.96155
$ws.Range("J22").Value = 899.8
$ws.Range("L22").Value = 899.8
$ws.Range("N22").Value = -1489.8
$ws.Range("H27").Value = 830.96155
$ws.Range("J27").Value = 899.8
$ws.Range("L27").Value = 899.8
$ws.Range("N27").Value = -1113.8
$ws.Range("H121").Value = 69997
$ws.Range("J121").Value = 69997
$ws.Range("L121").Value = 69997
$ws.Range("N121").Value = -73491
$ws.Range("H132").Value = 2202.7593
$ws.Range("I132").Value = 1899.14
$ws.Range("K132").Value = 5697.42
$ws.Range("M132").Value = -3167.42

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 29995
$ws.Range("J121").Value = 29995
$ws.Range("L121").Value = 29995
$ws.Range("N121").Value = -33489
$ws.Range("H122").Value = 1446.8214
$ws.Range("I122").Value = 1388.875
$ws.Range("K122").Value = 4166.625
$ws.Range("M122").Value = -1716.625
